# Performer Process Ongoing. Review Sheet logic was added and corrections in the main workflow.
#
# The "Constants" worksheet holds a Name/Value/Description settings table.
# This change:
#   - replaces the old PathTemplateReviewSheet setting (row 19) with a new
#     PathTemplatesFile setting pointing at a combined templates workbook
#   - inserts 5 new rows (20-24) adding PathMarginsList, PathCustomerNameList,
#     CustomerNameList_WorksheetName (+description) and
#     TemplateFile_WorksheetReviewSheetTemplate
#   - adds a new OtherDeductionsList_SheetClickOnText row after the other
#     OtherDeductionsList_Sheet* rows (now row 39)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Insert 5 new rows (20-24) for the new Path*/CustomerNameList settings --
$ws.Rows("20:24").Insert()
$ws.Rows("20:24").RowHeight = 14.25

$ws.Range("A20").Value = "PathMarginsList"
$ws.Range("B20").Value = "Data\Margins List.xlsx"

$ws.Range("A21").Value = "PathCustomerNameList"
$ws.Range("B21").Value = "Data\Customer Account List.xlsx"

# row 22 stays blank (separator row, matches the pattern used elsewhere)

$ws.Range("A23").Value = "CustomerNameList_WorksheetName"
$ws.Range("B23").Value = "Sheet1"
$ws.Range("C23").Value = "The sheet name of the customer name list file."

# --- Row 19: PathTemplateReviewSheet -> PathTemplatesFile -------------------
$ws.Range("A19").Value = "PathTemplatesFile"
$ws.Range("B19").Value = "Data\Templatesxlsx"

$ws.Range("A24").Value = "TemplateFile_WorksheetReviewSheetTemplate"
$ws.Range("B24").Value = "Review Sheet"

# --- New OtherDeductionsList_SheetClickOnText row (now row 39) -------------
$ws.Range("A39").Value = "OtherDeductionsList_SheetClickOnText"
$ws.Range("B39").Value = "ClickOnText"

# --- Update view/selection to match the saved state -------------------------
$ws.Activate() | Out-Null
$ws.Range("A39").Select() | Out-Null
